# Added "transaction finding by amount" test case.
#
# The old TC-009 ("Accounts overview" balance-check test) is removed, and the
# two test cases that followed it (TC-010 and TC-011, both covering the
# "Find Transactions" feature - one by a valid amount, one by an invalid
# amount) are shifted up to become the new TC-009 and TC-010.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old TC-009 row (row 16) entirely - clear its contents and reset
# the row height back to the sheet's standard height (no more custom height).
$ws.Range("A16:F16").ClearContents()
$ws.Rows.Item(16).AutoFit()

# The old TC-010 (row 18) becomes the new TC-009.
$ws.Range("C18").Value = "TC-009"

# The old TC-011 (row 19) becomes the new TC-010.
$ws.Range("C19").Value = "TC-010"

# Scroll the view up a bit now that a row's worth of content moved up
# (previously the view was scrolled so row 16 was at the top; now row 14).
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 2
